# This workbook is a KiBot/KiCost-generated Bill of Materials / Costs report.
# The commit simply re-ran the generator (new CI build), which bumped a
# handful of textual values that appear on every BoM-ish sheet:
#   - the revision string "r1"                         -> "${git_hash}"
#   - the generation date "2020-09-09"                  -> "${date}"
#   - the KiCad version    "6.0.10+dfsg-1~bpo11+1"       -> "6.0.11+dfsg-1~bpo11+1"
#   - the "Created:" timestamp on the Costs sheets       -> new timestamp
#   - the KiCost/KiBot version banner                    -> new KiBot point release
#
# These values are duplicated (as literal cells, backed by the shared
# string table) on the BoM, DNF, Costs and "Costs (DNF)" sheets, so update
# each one explicitly.

$wb = $excel.ActiveWorkbook

$revisionOld = "r1"
$revisionNew = "`${git_hash}"
$dateOld = "2020-09-09"
$dateNew = "`${date}"
$kicadVerOld = "6.0.10+dfsg-1~bpo11+1"
$kicadVerNew = "6.0.11+dfsg-1~bpo11+1"
$createdOld = "2023-02-06 16:47:21"
$createdNew = "2023-04-02 13:40:07"
$kicostOld = "KiCost® v1.1.15 + KiBot v1.6.0"
$kicostNew = "KiCost® v1.1.15 + KiBot v1.6.1"

# --- BoM + DNF sheets: Revision / Date / KiCad Version live in C4:D4, C5:D5, C6:D6 ---
foreach ($sheetName in @("BoM", "DNF")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("D4").Value = $revisionNew
    $ws.Range("D5").Value = $dateNew
    $ws.Range("D6").Value = $kicadVerNew
}

# --- Costs + Costs (DNF) sheets: same trio lives in D4:E4, D5:E5, D6:E6, plus
#     a "Created:" timestamp and the KiCost/KiBot banner further down ---
$ws = $wb.Worksheets.Item("Costs")
$ws.Range("E4").Value = $revisionNew
$ws.Range("E5").Value = $dateNew
$ws.Range("E6").Value = $kicadVerNew
$ws.Range("B24").Value = $createdNew
$ws.Range("A25").Value = $kicostNew

$ws = $wb.Worksheets.Item("Costs (DNF)")
$ws.Range("E4").Value = $revisionNew
$ws.Range("E5").Value = $dateNew
$ws.Range("E6").Value = $kicadVerNew
$ws.Range("B21").Value = $createdNew
$ws.Range("A22").Value = $kicostNew
